$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.908.66"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "2.306.64"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.71"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.19"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("E7").Value = "  -1.56%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -2.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.68"
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.39"
$ws.Range("E12").Value = "  +3.46%  "
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.77"
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "2.665.51"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "2.311.73"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.783"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "42.839.20"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.11"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.59"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.94"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.53"
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.42"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.09"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.03"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.01"
$ws.Range("E35").Value = "  -3.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.52"
$ws.Range("E36").Value = "  -3.71%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("D43").Value = "2.006.91"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").Value = "  -2.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.40"
$ws.Range("E45").Value = "  +4.52%  "
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("E47").Value = "  -3.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("E48").Value = "  -3.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.87"
$ws.Range("E49").Value = "  +5.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.97"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").Value = "2.533.64"
$ws.Range("E51").Value = "  -0.22%  "
